$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws3 = $wb.Worksheets.Item("2017 LEAVE BALANCE")
$ws4 = $wb.Worksheets.Item("CONVERTION")

# ---------------------------------------------------------------------------
# NOTE: shared-string table order matters - new strings are appended to the
# shared string table in the exact order in which they are first assigned.
# We therefore set the string-valued cells that introduce brand-new strings
# in the same order they appear in the target sharedStrings.xml.
# ---------------------------------------------------------------------------

# --- Sheet3 (2017 LEAVE BALANCE): rows 32-41 new data --------------------

$ws3.Range("A32").Value = 45200
$ws3.Range("B32").Value = "VL(4-0-0)"
$ws3.Range("D32").Value = 4
$ws3.Range("K32").Value = "10/12-13, 19-20/2023"          # new shared string 71

$ws3.Range("B33").Value = "SL(2-0-0)"
$ws3.Range("H33").Value = 2
$ws3.Range("K33").Value = "10/5-6/2023"                    # new shared string 72

# --- Sheet4 (CONVERTION): A6 label ---------------------------------------
$ws4.Range("A6").Value = "TOTAL LEAVE BALANCE"              # new shared string 73
$ws4.Range("D2").Copy()
$ws4.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("B34").Value = "VL(4-0-0)"
$ws3.Range("D34").Value = 4
$ws3.Range("K34").Value = "10/26-29/2023"                   # new shared string 74

$ws3.Range("A35").Value = 45231
$ws3.Range("B35").Value = "VL(6-0-0)"
$ws3.Range("D35").Value = 6
$ws3.Range("K35").Value = "11/9,10,16,17,23,24/2023"        # new shared string 75

$ws3.Range("B36").Value = "SL(1-0-0)"
$ws3.Range("H36").Value = 1
$ws3.Range("K11").Copy()
$ws3.Range("K36").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws3.Range("K36").Value = 45233

$ws3.Range("A37").Value = 45261
$ws3.Range("B37").Value = "SL(2-0-0)"
$ws3.Range("H37").Value = 2
$ws3.Range("K37").Value = "12/2,3/2023"                     # new shared string 76

$ws3.Range("B38").Value = "VL(5-0-0)"
$ws3.Range("D38").Value = 5
$ws3.Range("K38").Value = "12/1,6,7,14,15/2023"             # new shared string 77

$ws3.Range("B39").Value = "SL(2-0-0)"
$ws3.Range("H39").Value = 2
$ws3.Range("K39").Value = "11/29,30/2023"                   # new shared string 78

$ws3.Range("D40").Value = 10
$ws3.Range("K40").Value = "12/9,10,13-17,20-22/2023"        # new shared string 79

$ws3.Range("B40").Value = "VL(10-0-0)"                      # new shared string 80

# --- Sheet2 (2018 LEAVE CREDITS) year-end rows ---------------------------
$ws2.Range("C84").Value = 1.25
$ws2.Range("C85").Value = 1.25
$ws2.Range("C86").Value = 1.25

$ws2.Range("B87").Value = "FL(5-0-0)"
$ws2.Range("D87").Value = 5

$ws2.Range("A10").Copy()
$ws2.Range("A88").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("A88").Value = "'2024"                           # new shared string 81 (forced text)

$ws2.Range("K87").Value = "12/23,24,27,28,29/2023"          # new shared string 82

$ws3.Range("B41").Value = "VL(2-0-0)"
$ws3.Range("D41").Value = 2
$ws3.Range("K41").Value = "12/30,31/2023"                   # new shared string 83

# ---------------------------------------------------------------------------
# Remaining non-string-table-affecting edits
# ---------------------------------------------------------------------------

# Sheet2: shift the date column down by one row (rows 89-95)
$ws2.Range("A89").Value = 45292
$ws2.Range("A90").Value = 45323
$ws2.Range("A91").Value = 45352
$ws2.Range("A92").Value = 45383
$ws2.Range("A93").Value = 45413
$ws2.Range("A94").Value = 45444
$ws2.Range("A95").Value = 45474

# Sheet2: insert a new blank table row before the final (bottom-border) row.
# First clone the final row's (135) styling/formula down into the brand-new
# row 136, then restyle row 135 itself as a normal interior row.
$ws2.Range("A135:K135").Copy()
$ws2.Range("A136:K136").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("G136").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

$ws2.Range("A134:K134").Copy()
$ws2.Range("A135:K135").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("G135").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

$lo1 = $ws2.ListObjects.Item("Table1")
$lo1.Resize($ws2.Range("A8:K136"))

# Sheet4: A7 total-leave-balance formula
$ws2.Range("E11").Copy()
$ws4.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws4.Range("A7").Formula = "=SUM('2018 LEAVE CREDITS'!E9,'2018 LEAVE CREDITS'!I9)"

# Sheet3: mergeCells need to be re-emitted in ascending order
$mergeRanges = @("B2:C2","F2:G2","J2:K2","B3:C3","F3:G3","J3:K3","B4:C4","F4:G4","J4:K4","C7:F7","G7:J7")
foreach ($r in $mergeRanges) {
    $ws3.Range($r).UnMerge()
}
foreach ($r in $mergeRanges) {
    $ws3.Range($r).Merge()
}

$excel.CalculateFullRebuild()
